$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update metric label text (shared-string content changes) ---
$ws.Range("C64").Value = "Blocajes Exitosos Fuera Área Peq."
$ws.Range("C65").Value = "Blocajes Fallidos Fuera Área Peq."
$ws.Range("C66").Value = "% Éxito Blocajes Fuera Área Peq."
$ws.Range("C85").Value = "Errores que terminan en Tiro"

# --- Update numeric ponderacion values in column E ---
$ws.Range("E58").Value = 1
$ws.Range("E59").Value = 3
$ws.Range("E65").Value = 3
$ws.Range("E68").Value = 3
$ws.Range("E71").Value = 1
$ws.Range("E72").Value = 3

# --- Apply AutoFilter on categoria (column A) keeping only Construcción / Juego Aéreo ---
$ws.Range("A1:E85").AutoFilter(1, @("Construcción", "Juego Aéreo"), 7)

# --- Update the saved view state: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D85").Select()
